# Auto-generated Excel COM-interop script
# Applies refreshed market-price / profit values to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 253.75
$ws.Range("I8").Value = 224.5
$ws.Range("K8").Value = 673.5
$ws.Range("M8").Value = -534.5

# Row 39
$ws.Range("H39").Value = 3922.3704
$ws.Range("I39").Value = 221.125
$ws.Range("J39").Value = 33532.332
$ws.Range("K39").Value = 663.375
$ws.Range("L39").Value = 100596.996
$ws.Range("M39").Value = -367.375
$ws.Range("N39").Value = -101188.996

# Row 64
$ws.Range("H64").Value = 5931.857
$ws.Range("I64").Value = 4893
$ws.Range("J64").Value = 6347.4
$ws.Range("K64").Value = 4893
$ws.Range("L64").Value = 6347.4
$ws.Range("M64").Value = -4645
$ws.Range("N64").Value = -6843.4

# Row 67
$ws.Range("H67").Value = 5931.857
$ws.Range("I67").Value = 4893
$ws.Range("J67").Value = 6347.4
$ws.Range("K67").Value = 4893
$ws.Range("L67").Value = 6347.4
$ws.Range("M67").Value = -4035
$ws.Range("N67").Value = -8063.4

# Row 88
$ws.Range("H88").Value = 1219.75
$ws.Range("I88").Value = 1257.3334
$ws.Range("K88").Value = 1257.3334
$ws.Range("M88").Value = -851.3334

# Row 91
$ws.Range("H91").Value = 1219.75
$ws.Range("I91").Value = 1257.3334
$ws.Range("K91").Value = 1257.3334
$ws.Range("M91").Value = 146.6666

# Row 116
$ws.Range("H116").Value = 5735.316
$ws.Range("I116").Value = 5726.9287
$ws.Range("J116").Value = 5758.8
$ws.Range("K116").Value = 5726.9287
$ws.Range("L116").Value = 5758.8
$ws.Range("M116").Value = -2284.9287
$ws.Range("N116").Value = -12642.8

# Row 138
$ws.Range("H138").Value = 4328.655
$ws.Range("J138").Value = 4466.9424
$ws.Range("L138").Value = 13400.8272
$ws.Range("N138").Value = -23680.8272

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1466.3043
$ws.Range("J45").Value = 1580
$ws.Range("L45").Value = 1580
$ws.Range("N45").Value = -2334

# Row 61
$ws.Range("H61").Value = 90890.25
$ws.Range("I61").Value = 3377.6
$ws.Range("J61").Value = 528453.5
$ws.Range("K61").Value = 3377.6
$ws.Range("L61").Value = 528453.5
$ws.Range("M61").Value = -3165.6
$ws.Range("N61").Value = -528877.5

# Row 136
$ws.Range("H136").Value = 90890.25
$ws.Range("I136").Value = 3377.6
$ws.Range("J136").Value = 528453.5
$ws.Range("K136").Value = 10132.8
$ws.Range("L136").Value = 1585360.5
$ws.Range("M136").Value = -7582.799999999999
$ws.Range("N136").Value = -1590460.5

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 20475
$ws.Range("I22").Value = 20475
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 20475
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -20302
$ws.Range("N22").ClearContents()

# Row 86
$ws.Range("H86").Value = 2504.08
$ws.Range("I86").Value = 2480.4736
$ws.Range("J86").Value = 2578.8333
$ws.Range("K86").Value = 2480.4736
$ws.Range("L86").Value = 2578.8333
$ws.Range("M86").Value = -1357.4736
$ws.Range("N86").Value = -4824.8333

# Row 89
$ws.Range("H89").Value = 2504.08
$ws.Range("I89").Value = 2480.4736
$ws.Range("J89").Value = 2578.8333
$ws.Range("K89").Value = 12402.368
$ws.Range("L89").Value = 12894.1665
$ws.Range("M89").Value = -6786.367999999999
$ws.Range("N89").Value = -24126.1665

# Row 107
$ws.Range("H107").Value = 913.5714
$ws.Range("I107").Value = 1505.5
$ws.Range("K107").Value = 1505.5
$ws.Range("M107").Value = 414.5

# Row 134
$ws.Range("H134").Value = 59237.863
$ws.Range("I134").Value = 143656.58
$ws.Range("K134").Value = 430969.74
$ws.Range("M134").Value = -428434.74

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 552.1905
$ws.Range("I7").Value = 601.6316
$ws.Range("J7").Value = 82.5
$ws.Range("K7").Value = 601.6316
$ws.Range("L7").Value = 82.5
$ws.Range("M7").Value = -488.6316
$ws.Range("N7").Value = -308.5

# Row 31
$ws.Range("H31").Value = 7464.121
$ws.Range("I31").Value = 1205.12
$ws.Range("J31").Value = 27023.5
$ws.Range("K31").Value = 1205.12
$ws.Range("L31").Value = 27023.5
$ws.Range("M31").Value = -910.1199999999999
$ws.Range("N31").Value = -27613.5

# Row 34
$ws.Range("H34").Value = 7464.121
$ws.Range("I34").Value = 1205.12
$ws.Range("J34").Value = 27023.5
$ws.Range("K34").Value = 1205.12
$ws.Range("L34").Value = 27023.5
$ws.Range("M34").Value = -1003.12
$ws.Range("N34").Value = -27427.5

# Row 52
$ws.Range("H52").Value = 91000
$ws.Range("I52").Value = 32000
$ws.Range("K52").Value = 32000
$ws.Range("M52").Value = -31706

# Row 132
$ws.Range("H132").Value = 30305024
$ws.Range("I132").Value = 2099.3447
$ws.Range("K132").Value = 6298.034100000001
$ws.Range("M132").Value = -3768.034100000001

# Row 134
$ws.Range("H134").Value = 66676336
$ws.Range("I134").Value = 2698.3
$ws.Range("J134").Value = 200023620
$ws.Range("K134").Value = 8094.900000000001
$ws.Range("L134").Value = 600070860
$ws.Range("M134").Value = -5559.900000000001
$ws.Range("N134").Value = -600075930

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 4239.5
$ws.Range("J22").Value = 4239.5
$ws.Range("L22").Value = 12718.5
$ws.Range("N22").Value = -13056.5

# Row 27
$ws.Range("H27").Value = 4239.5
$ws.Range("J27").Value = 4239.5
$ws.Range("L27").Value = 12718.5
$ws.Range("N27").Value = -12922.5

# Row 34
$ws.Range("H34").Value = 3129.0833
$ws.Range("J34").Value = 5704.5
$ws.Range("L34").Value = 17113.5
$ws.Range("N34").Value = -17281.5

# Row 39
$ws.Range("H39").Value = 7470.857
$ws.Range("I39").Value = 7499.5
$ws.Range("J39").Value = 7459.4
$ws.Range("K39").Value = 22498.5
$ws.Range("L39").Value = 22378.2
$ws.Range("M39").Value = -22204.5
$ws.Range("N39").Value = -22966.2

# Row 55
$ws.Range("H55").Value = 2076.75
$ws.Range("I55").Value = 2328.6
$ws.Range("J55").Value = 1657
$ws.Range("K55").Value = 6985.799999999999
$ws.Range("L55").Value = 4971
$ws.Range("M55").Value = -6808.799999999999
$ws.Range("N55").Value = -5325

# Row 131
$ws.Range("H131").Value = 1493.91
$ws.Range("I131").Value = 1048.3334
$ws.Range("J131").Value = 1522.3511
$ws.Range("K131").Value = 3145.0002
$ws.Range("L131").Value = 4567.0533
$ws.Range("M131").Value = 1894.9998
$ws.Range("N131").Value = -14647.0533

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -388
$ws.Range("N5").ClearContents()

# Row 26
$ws.Range("H26").Value = 23000
$ws.Range("J26").Value = 23000
$ws.Range("L26").Value = 23000
$ws.Range("N26").Value = -23560

# Row 50
$ws.Range("H50").Value = 23000
$ws.Range("J50").Value = 23000
$ws.Range("L50").Value = 23000
$ws.Range("N50").Value = -23996

# Row 80
$ws.Range("H80").Value = 8147.5
$ws.Range("I80").Value = 8147.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 8147.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -7149.5
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 8147.5
$ws.Range("I83").Value = 8147.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 40737.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -35745.5
$ws.Range("N83").ClearContents()

# Row 132
$ws.Range("H132").Value = 27597.1
$ws.Range("I132").Value = 6376.6
$ws.Range("J132").Value = 48817.6
$ws.Range("K132").Value = 19129.8
$ws.Range("L132").Value = 146452.8
$ws.Range("M132").Value = -16599.8
$ws.Range("N132").Value = -151512.8

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1215.3636
$ws.Range("I22").Value = 1455.5714
$ws.Range("K22").Value = 1455.5714
$ws.Range("M22").Value = -1160.5714

# Row 26
$ws.Range("H26").Value = 36995
$ws.Range("I26").Value = 36995
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 36995
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -36700
$ws.Range("N26").ClearContents()

# Row 27
$ws.Range("H27").Value = 1215.3636
$ws.Range("I27").Value = 1455.5714
$ws.Range("K27").Value = 1455.5714
$ws.Range("M27").Value = -1348.5714

# Row 40
$ws.Range("H40").Value = 5494.8335
$ws.Range("I40").Value = 4991.3335
$ws.Range("J40").Value = 5998.3335
$ws.Range("K40").Value = 4991.3335
$ws.Range("L40").Value = 5998.3335
$ws.Range("M40").Value = -4855.3335
$ws.Range("N40").Value = -6270.3335

# Row 82
$ws.Range("H82").Value = 1397.625
$ws.Range("I82").Value = 1906
$ws.Range("K82").Value = 1906
$ws.Range("M82").Value = -1545

# Row 85
$ws.Range("H85").Value = 1397.625
$ws.Range("I85").Value = 1906
$ws.Range("K85").Value = 1906
$ws.Range("M85").Value = -658

# Row 100
$ws.Range("H100").Value = 3182.1667
$ws.Range("J100").Value = 3448.1667
$ws.Range("L100").Value = 3448.1667
$ws.Range("N100").Value = -4530.1667

# Row 112
$ws.Range("H112").Value = 100387
$ws.Range("J112").Value = 100387
$ws.Range("L112").Value = 100387
$ws.Range("N112").Value = -103341

# Row 136
$ws.Range("H136").Value = 211474.73
$ws.Range("I136").Value = 25135.334
$ws.Range("K136").Value = 75406.00199999999
$ws.Range("M136").Value = -72856.00199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 15340.611
$ws.Range("J62").Value = 14888.111
$ws.Range("L62").Value = 14888.111
$ws.Range("N62").Value = -16136.111

# Row 65
$ws.Range("H65").Value = 15340.611
$ws.Range("J65").Value = 14888.111
$ws.Range("L65").Value = 74440.55500000001
$ws.Range("N65").Value = -80680.55500000001

# Row 96
$ws.Range("H96").Value = 1194.6364
$ws.Range("I96").Value = 948.6667
$ws.Range("J96").Value = 1489.8
$ws.Range("K96").Value = 948.6667
$ws.Range("L96").Value = 1489.8
$ws.Range("M96").Value = 424.3333
$ws.Range("N96").Value = -4235.8

# Row 103
$ws.Range("H103").Value = 24722
$ws.Range("J103").Value = 24722
$ws.Range("L103").Value = 24722
$ws.Range("N103").Value = -27066

# Row 107
$ws.Range("H107").Value = 1142.84
$ws.Range("I107").Value = 1224.5238
$ws.Range("J107").Value = 714
$ws.Range("K107").Value = 3673.5714
$ws.Range("L107").Value = 2142
$ws.Range("M107").Value = -1753.5714
$ws.Range("N107").Value = -5982

# Row 138
$ws.Range("H138").Value = 100319
$ws.Range("J138").Value = 100319
$ws.Range("L138").Value = 100319
$ws.Range("N138").Value = -110599
